$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, new Price (D) value ($null = leave unchanged), new Volume(1h) (E) value
$updates = @(
    @(2, "61.733.08", "  -1.85%  "),
    @(3, "3.005.55", "  -1.73%  "),
    @(4, $null, "  -0.03%  "),
    @(5, "543.63", "  +1.01%  "),
    @(6, $null, "  -3.09%  "),
    @(7, $null, "  +0.05%  "),
    @(8, "2.998.65", "  -1.75%  "),
    @(9, $null, "  -0.31%  "),
    @(10, "6.15", "  -1.63%  "),
    @(11, $null, "  -5.84%  "),
    @(12, $null, "  -1.42%  "),
    @(13, "34.52", "  +0.60%  "),
    @(14, $null, "  -0.87%  "),
    @(15, "3.493.98", "  -1.80%  "),
    @(16, "61.761.57", "  -1.84%  "),
    @(17, $null, "  -2.63%  "),
    @(18, "3.004.75", "  -1.88%  "),
    @(19, $null, "  +0.37%  "),
    @(20, "482.37", "  +3.17%  "),
    @(21, "13.23", "  -1.94%  "),
    @(22, $null, "  -3.40%  "),
    @(23, "6.95", "  -0.72%  "),
    @(24, "81.90", "  +4.63%  "),
    @(25, $null, "  -0.88%  "),
    @(26, $null, "  +0.06%  "),
    @(27, $null, "  +0.35%  "),
    @(28, "7.72", "  -1.62%  "),
    @(29, $null, "  -0.21%  "),
    @(30, $null, "  +2.23%  "),
    @(31, "25.65", "  -1.32%  "),
    @(32, "1.12", "  -2.84%  "),
    @(33, "5.65", "  +4.00%  "),
    @(34, $null, "  +0.78%  "),
    @(35, "55.05", "  -6.10%  "),
    @(36, $null, "  -1.47%  "),
    @(37, "443.44", "  -7.57%  "),
    @(38, "3.141.56", "  -3.09%  "),
    @(39, $null, "  +0.83%  "),
    @(40, $null, "  -3.60%  "),
    @(41, "0.118", "  +0.30%  "),
    @(42, "8.07", "  -0.46%  "),
    @(43, $null, "  -5.08%  "),
    @(44, "26.42", "  +5.57%  "),
    @(45, $null, "  -0.06%  "),
    @(46, $null, "  -3.20%  "),
    @(47, $null, "  +0.13%  "),
    @(48, $null, "  -2.35%  "),
    @(49, "115.34", "  -6.08%  "),
    @(50, "1.29", "  +5.10%  "),
    @(51, $null, "  -5.92%  ")
)

# Temporarily mark the Price/Volume columns as Text so that numeric-looking
# strings (e.g. "6.15") are written verbatim instead of being coerced into
# Excel numbers, then clear the formatting again so cell styling is left
# exactly as it was before the edit.
$updateRange = $ws.Range("D2:E51")
$updateRange.NumberFormat = "@"

foreach ($u in $updates) {
    $row = $u[0]
    $dVal = $u[1]
    $eVal = $u[2]

    if ($null -ne $dVal) {
        $ws.Cells.Item($row, 4).Value = $dVal
    }
    $ws.Cells.Item($row, 5).Value = $eVal
}

$updateRange.ClearFormats()

